# Add new todo items (order chosen to reproduce shared-string insertion order)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C13").Value = "CRON JOB"
$ws.Range("A9").Value = "HATA SISTEMI"
$ws.Range("A10").Value = "POSITON BULUP GUNCELLEME ( CACHE ALININCA :)"
$ws.Range("A11").Value = "tarih null dönme hatası ( boş yap )"
$ws.Range("C14").Value = "Önerilenlere id verme"
$ws.Range("A13").Value = "Girişte banner ( splash screen )"
$ws.Range("A16").Value = "KURUM MOIL APP MAKALE SİTE İÇİN"
$ws.Range("A8").Value = "Slider"

# Mark completed items with strikethrough (hata azaltıldı)
$ws.Range("A2").Font.Strikethrough = $true
$ws.Range("A4").Font.Strikethrough = $true
$ws.Range("A5").Font.Strikethrough = $true
$ws.Range("A13").Font.Strikethrough = $true

# Update selection / page setup
$ws.Range("A9").Select() | Out-Null
$ws.PageSetup.Orientation = 1
